$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Delete" test-suite table (columns A:E, rows 2..) gets a new row
# inserted at row 4 for the "JobHistoryErrors" automateable tag, shifting the
# existing rows 4 and 5 down to rows 5 and 6. The summary box in columns
# F:I is not part of this table and stays where it is; its formulas are
# whole-column (SUM($B:$B), SUM($C:$C), COUNTIF($D:$D,...), etc.) so they
# simply recalculate once the underlying data moves.
# ---------------------------------------------------------------------------

# Push old row 5 (values + formatting) down into row 6, then old row 4 down
# into row 5 -- bottom-up so we never clobber a row before it has been
# copied elsewhere.
$ws.Range("A5:E5").Copy($ws.Range("A6:E6"))
$ws.Range("A4:E4").Copy($ws.Range("A5:E5"))

# Write the new "JobHistoryErrors" row into the now-vacated row 4.
$ws.Range("A4").Value = "JobHistoryErrors"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "Finished"
$ws.Range("E4").ClearContents()

# Move the active selection the way the author left it.
$ws.Range("A8").Select()
